$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2302.3845
$ws.Range("J38").Value = 3626.125
$ws.Range("L38").Value = 10878.375
$ws.Range("N38").Value = -11622.375
$ws.Range("H43").Value = 1128.8572
$ws.Range("I43").Value = 1025.5
$ws.Range("J43").Value = 1266.6666
$ws.Range("K43").Value = 1025.5
$ws.Range("L43").Value = 1266.6666
$ws.Range("M43").Value = -956.5
$ws.Range("N43").Value = -1404.6666
$ws.Range("H62").Value = 4286435.5
$ws.Range("I62").Value = 6545107
$ws.Range("J62").Value = 20055.555
$ws.Range("K62").Value = 6545107
$ws.Range("L62").Value = 20055.555
$ws.Range("M62").Value = -6544483
$ws.Range("N62").Value = -21303.555
$ws.Range("H65").Value = 4286435.5
$ws.Range("I65").Value = 6545107
$ws.Range("J65").Value = 20055.555
$ws.Range("K65").Value = 32725535
$ws.Range("L65").Value = 100277.775
$ws.Range("M65").Value = -32722415
$ws.Range("N65").Value = -106517.775
$ws.Range("H82").Value = 1708.7778
$ws.Range("I82").Value = 1354.1428
$ws.Range("J82").Value = 2950
$ws.Range("K82").Value = 4062.4284
$ws.Range("L82").Value = 8850
$ws.Range("M82").Value = -3656.4284
$ws.Range("N82").Value = -9662
$ws.Range("H85").Value = 1708.7778
$ws.Range("I85").Value = 1354.1428
$ws.Range("J85").Value = 2950
$ws.Range("K85").Value = 4062.4284
$ws.Range("L85").Value = 8850
$ws.Range("M85").Value = -2658.4284
$ws.Range("N85").Value = -11658
$ws.Range("H137").Value = 21740228
$ws.Range("I137").Value = 32258868
$ws.Range("J137").Value = 1704.0667
$ws.Range("K137").Value = 96776604
$ws.Range("L137").Value = 5112.2001
$ws.Range("M137").Value = -96774054
$ws.Range("N137").Value = -10212.2001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15191.0625
$ws.Range("I32").Value = 2629.923
$ws.Range("K32").Value = 2629.923
$ws.Range("M32").Value = -2342.923
$ws.Range("H61").Value = 1394.1964
$ws.Range("I61").Value = 1093.1346
$ws.Range("J61").Value = 5308
$ws.Range("K61").Value = 1093.1346
$ws.Range("L61").Value = 5308
$ws.Range("M61").Value = -881.1346000000001
$ws.Range("N61").Value = -5732
$ws.Range("H63").Value = 12908.909
$ws.Range("I63").Value = 22383
$ws.Range("J63").Value = 1540
$ws.Range("K63").Value = 22383
$ws.Range("L63").Value = 1540
$ws.Range("M63").Value = -21697
$ws.Range("N63").Value = -2912
$ws.Range("H66").Value = 12908.909
$ws.Range("I66").Value = 22383
$ws.Range("J66").Value = 1540
$ws.Range("K66").Value = 111915
$ws.Range("L66").Value = 7700
$ws.Range("M66").Value = -108483
$ws.Range("N66").Value = -14564
$ws.Range("H74").Value = 3793.7114
$ws.Range("I74").Value = 1262.6216
$ws.Range("K74").Value = 1262.6216
$ws.Range("M74").Value = -388.6215999999999
$ws.Range("H77").Value = 3793.7114
$ws.Range("I77").Value = 1262.6216
$ws.Range("K77").Value = 6313.108
$ws.Range("M77").Value = -1945.108
$ws.Range("H122").Value = 9342.267
$ws.Range("I122").Value = 10961.167
$ws.Range("K122").Value = 32883.501
$ws.Range("M122").Value = -30433.501
$ws.Range("H136").Value = 1394.1964
$ws.Range("I136").Value = 1093.1346
$ws.Range("J136").Value = 5308
$ws.Range("K136").Value = 3279.4038
$ws.Range("L136").Value = 15924
$ws.Range("M136").Value = -729.4038
$ws.Range("N136").Value = -21024

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1714.5
$ws.Range("I31").Value = 983.69446
$ws.Range("J31").Value = 4345.4
$ws.Range("K31").Value = 983.69446
$ws.Range("L31").Value = 4345.4
$ws.Range("M31").Value = -688.69446
$ws.Range("N31").Value = -4935.4
$ws.Range("H34").Value = 1714.5
$ws.Range("I34").Value = 983.69446
$ws.Range("J34").Value = 4345.4
$ws.Range("K34").Value = 983.69446
$ws.Range("L34").Value = 4345.4
$ws.Range("M34").Value = -781.69446
$ws.Range("N34").Value = -4749.4
$ws.Range("J86").Value = 1692.1875
$ws.Range("L86").Value = 1692.1875
$ws.Range("N86").Value = -3938.1875
$ws.Range("J89").Value = 1692.1875
$ws.Range("L89").Value = 8460.9375
$ws.Range("N89").Value = -19692.9375
$ws.Range("H99").Value = 20833608
$ws.Range("I99").Value = 20833608
$ws.Range("K99").Value = 20833608
$ws.Range("M99").Value = -20832110
$ws.Range("H126").Value = 20833608
$ws.Range("I126").Value = 20833608
$ws.Range("K126").Value = 62500824
$ws.Range("M126").Value = -62498354

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7894.59
$ws.Range("I4").Value = 111.117645
$ws.Range("J4").Value = 13909.091
$ws.Range("K4").Value = 333.352935
$ws.Range("L4").Value = 41727.273
$ws.Range("M4").Value = -221.352935
$ws.Range("N4").Value = -41951.273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1000
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("N5").Value = -1224
$ws.Range("H102").Value = 2697.2104
$ws.Range("I102").Value = 2487.25
$ws.Range("K102").Value = 2487.25
$ws.Range("M102").Value = -865.25
$ws.Range("H126").Value = 1740.3636
$ws.Range("I126").Value = 1106.1333
$ws.Range("J126").Value = 2268.889
$ws.Range("K126").Value = 3318.3999
$ws.Range("L126").Value = 6806.667
$ws.Range("M126").Value = -848.3998999999999
$ws.Range("N126").Value = -11746.667
$ws.Range("H132").Value = 3825.2122
$ws.Range("I132").Value = 3591.6155
$ws.Range("J132").Value = 4692.857
$ws.Range("K132").Value = 10774.8465
$ws.Range("L132").Value = 14078.571
$ws.Range("M132").Value = -8244.8465
$ws.Range("N132").Value = -19138.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4546829
$ws.Range("I16").Value = 5001464
$ws.Range("J16").Value = 478.5
$ws.Range("K16").Value = 5001464
$ws.Range("L16").Value = 478.5
$ws.Range("M16").Value = -5001294
$ws.Range("N16").Value = -818.5
$ws.Range("H40").Value = 3471.963
$ws.Range("I40").Value = 1940.8889
$ws.Range("J40").Value = 4237.5
$ws.Range("K40").Value = 1940.8889
$ws.Range("L40").Value = 4237.5
$ws.Range("M40").Value = -1804.8889
$ws.Range("N40").Value = -4509.5
$ws.Range("H100").Value = 2843859
$ws.Range("I100").Value = 10418817
$ws.Range("J100").Value = 3250
$ws.Range("K100").Value = 10418817
$ws.Range("L100").Value = 3250
$ws.Range("M100").Value = -10418276
$ws.Range("N100").Value = -4332
$ws.Range("H122").Value = 3390.25
$ws.Range("I122").Value = 1975
$ws.Range("K122").Value = 5925
$ws.Range("M122").Value = -3475

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 500001000
$ws.Range("I96").Value = 1000000000
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 1000000000
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = -999998627
$ws.Range("N96").Value = -4746
$ws.Range("H107").Value = 2646463.2
$ws.Range("I107").Value = 4274287
$ws.Range("K107").Value = 12822861
$ws.Range("M107").Value = -12820941
$ws.Range("H132").Value = 11365942
$ws.Range("I132").Value = 15627357
$ws.Range("J132").Value = 2166.6667
$ws.Range("K132").Value = 46882071
$ws.Range("L132").Value = 6500.000100000001
$ws.Range("M132").Value = -46879541
$ws.Range("N132").Value = -11560.0001
